$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("About")

# Update the note text to reference GST
$ws1.Range("A10").Value = "We are using the national average sales tax rate - GST"

# Move the selection to A11 as recorded in the saved file
$ws1.Range("A11").Select()
